# Replace the temporary/sample crew member in row 3 with the real crew
# member's details, and clear out the placeholder second crew row (row 4)
# that was only there for the import preview, per:
# "Modify Crew import from Excel (Add new Temporary Crew and the then Insert)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crew member in row 3: Name, Rank, Email
$ws.Range("B3").Value = "Bhuban"
$ws.Range("C3").Value = "Master"
$ws.Range("D3").Value = "master@gmail.com"

# Remove the temporary second crew row's data (Name/Rank/Email), leaving
# the numbering and formatting intact
$ws.Range("B4:D4").ClearContents()

# Leave selection on the edited email cell
$ws.Range("D3").Select()
